$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304"
for ($i = 1; $i -le 21; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $val = [string]$cell.Value2
    if ($val.EndsWith("_old")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2210"
    } elseif ($val.EndsWith("_new")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
    }
}

# 2. Turn the used range into an Excel Table (ListObject)
$tableRange = $ws.Range("A1:U68")
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"

# 3. Freeze panes at row 2 (split below header row)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
